$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ORG_NAT_IDENOLD"
$ws.Range("E1").Value = "ORG_NAT_IDENNEW"
$ws.Range("F1").Value = "ORG_NAT_STATUS"

$ws.Range("D1:F1").HorizontalAlignment = -4131

$ws.Range("F2").Select()
